$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E11").Value = 229286
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 4

$ws.Rows("19:22").Delete()

$ws.Range("C16").Value = "1002187459"
$ws.Range("D16").Value = "ESTEFANIA GIL HINCAPIE"
$ws.Range("E16").Value = "2502"
$ws.Range("F16").Value = 13286
$ws.Range("G16").Value = 1600000

$ws.Range("C17").Value = "1044930543"
$ws.Range("D17").Value = "YORDI JAVIER GOMEZ MAZA"
$ws.Range("E17").Value = "2506"
$ws.Range("F17").Value = 72000
$ws.Range("G17").Value = 1800000

$ws.Range("C18").Value = "1044930543"
$ws.Range("D18").Value = "YORDI JAVIER GOMEZ MAZA"
$ws.Range("E18").Value = "2507"
$ws.Range("F18").Value = 72000
$ws.Range("G18").Value = 1800000

$ws.Range("C19").Value = "1044930543"
$ws.Range("D19").Value = "YORDI JAVIER GOMEZ MAZA"
$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 72000
$ws.Range("G19").Value = 1800000
